# Applies the cryptos.xlsx price/volume refresh + three coin-row swaps
# described by the diff (rows 2-51), matching the commit's GitHub Actions
# scheduled data update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '63.565.77'
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -1.20%  '
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '3.075.20'
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -0.90%  '
$ws.Cells.Item(4, 5).Value = '  -0.73%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '592.77'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.81%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '154.88'
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +2.17%  '
$ws.Cells.Item(7, 5).Value = '  -0.30%  '
$ws.Cells.Item(8, 5).Value = '  +1.22%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '3.073.24'
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -0.83%  '
$ws.Cells.Item(10, 5).Value = '  -0.95%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '5.93'
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.25%  '
$ws.Cells.Item(12, 5).Value = '  -1.69%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '0.0000237'
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -2.07%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '36.69'
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -2.70%  '
$ws.Cells.Item(15, 5).Value = '  +1.19%  '
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '3.585.98'
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.88%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '7.19'
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '63.542.57'
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.52%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '3.077.52'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.99%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '481.88'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +3.00%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '14.48'
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -2.14%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '0.709'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -3.56%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '7.54'
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.17%  '
$ws.Cells.Item(24, 5).Value = '  +1.58%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '81.68'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.09%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '12.84'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -3.02%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '10.73'
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +10.56%  '
$ws.Cells.Item(28, 5).Value = '  +0.24%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '7.65'
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +3.27%  '
$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '2.69'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -0.18%  '
$ws.Cells.Item(31, 2).Value = 'ImmutableX'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +1.70%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.79%  '
$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '27.23'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -0.63%  '
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '0.112'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -3.22%  '
$ws.Cells.Item(35, 2).Value = 'Mantle'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '1.07'
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +1.19%  '
$ws.Cells.Item(36, 2).Value = 'PEPE'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '0.0₃0826'
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -1.63%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '6.06'
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -1.13%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '3.28'
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -1.17%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '2.22'
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -1.51%  '
$ws.Cells.Item(40, 5).Value = '  -1.18%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '50.61'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -0.17%  '
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '444.54'
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -1.33%  '
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '0.291'
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.50%  '
$ws.Cells.Item(44, 5).Value = '  +3.73%  '
$ws.Cells.Item(45, 5).Value = '  -1.64%  '
$ws.Cells.Item(46, 2).Value = 'Maker'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '2.825.32'
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -0.64%  '
$ws.Cells.Item(47, 2).Value = 'Arweave'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '39.83'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +4.13%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '132.31'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +1.86%  '
$ws.Cells.Item(49, 2).Value = 'USDe'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +0.02%  '
$ws.Cells.Item(50, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '25.23'
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.65%  '
$ws.Cells.Item(51, 5).Value = '  -1.03%  '
